$d = $word.ActiveDocument

$d.Content.Find.Execute("100÷3=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "799÷8=99, 7", 2) | Out-Null
$d.Content.Find.Execute("588÷2=294, 0", $true, $false, $false, $false, $false, $true, 1, $false, "521÷4=130, 1", 2) | Out-Null
$d.Content.Find.Execute("171÷3=57, 0", $true, $false, $false, $false, $false, $true, 1, $false, "767÷4=191, 3", 2) | Out-Null
$d.Content.Find.Execute("118÷9=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "904÷9=100, 4", 2) | Out-Null
$d.Content.Find.Execute("107÷7=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "726÷9=80, 6", 2) | Out-Null
$d.Content.Find.Execute("781÷5=156, 1", $true, $false, $false, $false, $false, $true, 1, $false, "661÷4=165, 1", 2) | Out-Null
$d.Content.Find.Execute("219÷5=43, 4", $true, $false, $false, $false, $false, $true, 1, $false, "995÷3=331, 2", 2) | Out-Null
$d.Content.Find.Execute("232÷9=25, 7", $true, $false, $false, $false, $false, $true, 1, $false, "123÷4=30, 3", 2) | Out-Null
$d.Content.Find.Execute("567÷3=189, 0", $true, $false, $false, $false, $false, $true, 1, $false, "853÷4=213, 1", 2) | Out-Null
$d.Content.Find.Execute("308÷4=77, 0", $true, $false, $false, $false, $false, $true, 1, $false, "448÷3=149, 1", 2) | Out-Null
$d.Content.Find.Execute("779÷9=86, 5", $true, $false, $false, $false, $false, $true, 1, $false, "303÷3=101, 0", 2) | Out-Null
$d.Content.Find.Execute("653÷7=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "208÷4=52, 0", 2) | Out-Null
$d.Content.Find.Execute("446÷2=223, 0", $true, $false, $false, $false, $false, $true, 1, $false, "426÷9=47, 3", 2) | Out-Null
$d.Content.Find.Execute("731÷6=121, 5", $true, $false, $false, $false, $false, $true, 1, $false, "132÷4=33, 0", 2) | Out-Null
$d.Content.Find.Execute("614÷9=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "401÷8=50, 1", 2) | Out-Null
$d.Content.Find.Execute("586÷8=73, 2", $true, $false, $false, $false, $false, $true, 1, $false, "490÷9=54, 4", 2) | Out-Null
$d.Content.Find.Execute("174÷8=21, 6", $true, $false, $false, $false, $false, $true, 1, $false, "273÷9=30, 3", 2) | Out-Null
$d.Content.Find.Execute("661÷5=132, 1", $true, $false, $false, $false, $false, $true, 1, $false, "859÷7=122, 5", 2) | Out-Null
$d.Content.Find.Execute("759÷4=189, 3", $true, $false, $false, $false, $false, $true, 1, $false, "467÷8=58, 3", 2) | Out-Null
$d.Content.Find.Execute("507÷8=63, 3", $true, $false, $false, $false, $false, $true, 1, $false, "854÷5=170, 4", 2) | Out-Null
$d.Content.Find.Execute("308÷2=154, 0", $true, $false, $false, $false, $false, $true, 1, $false, "731÷2=365, 1", 2) | Out-Null
$d.Content.Find.Execute("572÷7=81, 5", $true, $false, $false, $false, $false, $true, 1, $false, "941÷3=313, 2", 2) | Out-Null
$d.Content.Find.Execute("759÷6=126, 3", $true, $false, $false, $false, $false, $true, 1, $false, "961÷4=240, 1", 2) | Out-Null
$d.Content.Find.Execute("605÷9=67, 2", $true, $false, $false, $false, $false, $true, 1, $false, "380÷4=95, 0", 2) | Out-Null
$d.Content.Find.Execute("120÷7=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "698÷7=99, 5", 2) | Out-Null
